# Ontologia_DRONE.xlsx — "Add files via upload" re-edit
#
# The FatosIn sheet's row numbering column (A) had been left showing a
# stale "2" in every row; fix it so each row's index cell mirrors its own
# sheet row (A3=3, A4=4, ... A10=10 — A2 already read "2" and is untouched).
# Columns B:E are re-sized (closest to the Excel-autofit widths recorded in
# the saved file) and the live selection is moved to C2. NOW() in B18 on
# the Projeto sheet is volatile (ca="1") and is refreshed by recalculation
# automatically, so it needs no explicit write here.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FatosIn")
$ws.Activate()

# --- Row-index column A: rows 3-10 should hold their own row number ------
$ws.Range("A3").Value  = 3
$ws.Range("A4").Value  = 4
$ws.Range("A5").Value  = 5
$ws.Range("A6").Value  = 6
$ws.Range("A7").Value  = 7
$ws.Range("A8").Value  = 8
$ws.Range("A9").Value  = 9
$ws.Range("A10").Value = 10

# --- Column widths B:E (re-fit to content) --------------------------------
$ws.Columns.Item(2).ColumnWidth = 4.165    # B -> ~5
$ws.Columns.Item(3).ColumnWidth = 5.165    # C -> ~6.07421875
$ws.Columns.Item(4).ColumnWidth = 2.5      # D -> ~3.3828125
$ws.Columns.Item(5).ColumnWidth = 15.665   # E -> ~16.4609375

# --- Move the live selection to C2 ----------------------------------------
$null = $ws.Range("C2").Select()
